# Rerun experiment for xdribble model only
# Updates the "Total Minutes" (R-Proposed, column D) values on the three
# "Top 10 players ..." sheets and moves the active-cell selection on each.

$wb = $excel.ActiveWorkbook

# --- Sheet: Top 10 players goal 90 ---
$ws = $wb.Worksheets.Item("Top 10 players goal 90")
$ws.Range("D2").Value = 63
$ws.Range("D3").Value = 139
$ws.Range("D4").Value = 46
$ws.Range("D5").Value = 123
$ws.Range("D6").Value = 122
$ws.Range("D7").Value = 26
$ws.Range("D8").Value = 42
$ws.Range("D10").Value = 112
$ws.Range("D11").Value = 125
$ws.Activate()
$ws.Range("E11").Select()

# --- Sheet: Top 10 players assist 90 ---
$ws = $wb.Worksheets.Item("Top 10 players assist 90")
$ws.Range("D2").Value = 26
$ws.Range("D4").Value = 45
$ws.Range("D5").Value = 139
$ws.Range("D6").Value = 48
$ws.Range("D7").Value = 106
$ws.Range("D8").Value = 91
$ws.Range("D9").Value = 109
$ws.Range("D10").Value = 42
$ws.Range("D11").Value = 123
$ws.Activate()
$ws.Range("D10").Select()

# --- Sheet: Top 10 players goal assist 90 ---
$ws = $wb.Worksheets.Item("Top 10 players goal assist 90")
$ws.Range("D2").Value = 26
$ws.Range("D3").Value = 139
$ws.Range("D4").Value = 63
$ws.Range("D5").Value = 46
$ws.Range("D6").Value = 45
$ws.Range("D8").Value = 123
$ws.Range("D9").Value = 122
$ws.Range("D10").Value = 42
$ws.Range("D11").Value = 48
$ws.Activate()
$ws.Range("D10").Select()
